# Weekly data update: insert a new record (row) for this market/product
# right after the existing header + first data block, pushing all the
# historical rows (190-255) down by one, and append the old last row's
# data as the new row 256 (Excel's native "insert row" shifts everything
# down, which is exactly what happened in the source workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 190; this shifts rows 190:255 down to 191:256
# and keeps all existing formatting/row styles intact, just like Excel's
# Home > Insert > Insert Sheet Rows.
$ws.Rows.Item(190).Insert()

# Populate the newly inserted row 190 with this week's new record.
$ws.Cells.Item(190, 1).Value  = 4
$ws.Cells.Item(190, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(190, 3).Value  = "Los Lagos"
$ws.Cells.Item(190, 4).Value  = 44663
$ws.Cells.Item(190, 5).Value  = 10
$ws.Cells.Item(190, 6).Value  = 100112037
$ws.Cells.Item(190, 7).Value  = "Cebollín"
$ws.Cells.Item(190, 8).Value  = "Sin especificar"
$ws.Cells.Item(190, 9).Value  = "Primera"
$ws.Cells.Item(190, 10).Value = 140
$ws.Cells.Item(190, 11).Value = 9500
$ws.Cells.Item(190, 12).Value = 9500
$ws.Cells.Item(190, 13).Value = 9500
$ws.Cells.Item(190, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(190, 15).Value = "Región Metropolitana"
$ws.Cells.Item(190, 16).Value = 264
$ws.Cells.Item(190, 17).Value = 36
$ws.Cells.Item(190, 18).Value = "Hortaliza"
